$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) column values
$ws.Range("D2").Value = "26.516.28"
$ws.Range("D3").Value = "1.728.68"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.30"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4813"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2673"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06195"
$ws.Range("D10").Value = "1.730.07"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07189"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.60"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6105"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.531"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.17"
$ws.Range("D17").Value = "26.511.99"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9997"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006944"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.55"
$ws.Range("D21").Value = "1.953.71"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.525"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.817"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.258"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.95"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.771"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08037"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.693"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04522"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.618"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.001"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.081"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9102"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.373"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01503"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.45"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.553"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3878"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.979"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.794"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3402"

# Update Volume(1h) (E) column values
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  +3.02%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  +8.04%  "
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("E42").Value = "  -10.86%  "
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("E45").Value = "  +9.78%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("E51").Value = "  +0.38%  "
